$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple 1:1 text corrections (typos / wording fixes) in shared strings ---
$ws.Range("I2").Value = "어려운 파트에서 시간을 많이 투자해야 한다는 점이 아쉬움"
$ws.Range("N2").Value = "컨트롤 재밌다"
$ws.Range("H3").Value = "스토리가 좋은 게임을 좋아하는데, 어드벤처가 그나마 스토리도 좋고 플레이어가 캐릭터를 움직이거나 문제를 해결하는 식으로 게임이 진행되기 때문에 오래할 수 있음."
$ws.Range("T3").Value = "피지컬 싸움이 재밌음"
$ws.Range("P7").Value = "매 게임마다 수많은 변수들이 생기며 예상치 못한 문제나 사건들을 대응하는 것이 재밌음"
$ws.Range("X7").Value = "실제 플레이하는 것도 재미있고, 스스로의 팀을 꾸려나가는 재미를 많이 준다. 특히 fc 모바일 같은 게임의 경우 실제 선수들을 기반으로 하다보니 좋아하는 선수를 영입했을 때의 성취감이 더욱 크다."
$ws.Range("D10").Value = "스트레스가 풀려서 좋아요"
$ws.Range("J10").Value = "캐릭터 세지면 기분 좋음"
$ws.Range("D11").Value = "재밌다"
$ws.Range("D12").Value = "그냥 해보고 싶은 마음에 했는데 재미있어서 좋아하게 됐고"
$ws.Range("D14").Value = "총쏘기 시원하고 재밌게 때문입니다"
$ws.Range("R14").Value = "업데이트 할 때마다 재밌는 버전으로 나와서 재밌고"
$ws.Range("R15").Value = "다음 업데이트를 예측할 수 없어서 더 재밌게 하는거 같아요"
$ws.Range("D17").Value = "타격감 있다"
$ws.Range("R17").Value = "어린나이에 게임을 시작했고 배틀그라운드 게임이 제일 쉽다고 느껴짐"

# --- Row 6 rework: D6 typo fix, R6 split into R6 (truncated) + new S6, T6 replaced ---
$ws.Range("D6").Value = "1인칭 시점으로 플레이 하기 때문에 생동감을 느낄 수 있음. 특히 총으로 상대팀을 조준해서 죽일 때 가장 재미있음."
$ws.Range("R6").Value = "시간이 남을때 시간때우기 좋음 게임. 꾸준히 해야 게임 계정이 성장할 수 있는 점, 그리고 RPG 게임보다 사용자의 개입이 있고, 주로 두뇌를 써야되는 (몇 수를 보고 게임을 해야하는것) 게임이 많음. "
$ws.Range("S6").Value = "두뇌를 써야 하는 게임이기 때문에 게임에 질때, 스트레스를 유발할 수 있음"
$ws.Range("A6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("T6").Value = "1인칭 시점으로 플레이 하기 때문에 생동감을 느낄 수 있음. 특히 총으로 상대팀을 조준해서 죽일 때 가장 재미있음."

# --- Column width adjustments (D:AA), matching target widths as closely as the
#     pixel-quantized ColumnWidth setter allows (snaps to nearest 1/6 character unit) ---
$ws.Columns("D").ColumnWidth = 105.04666666666667
$ws.Columns("E:F").ColumnWidth = 61.29666666666667
$ws.Columns("G").ColumnWidth = 28.416666666666668
$ws.Columns("H").ColumnWidth = 149.41666666666666
$ws.Columns("I").ColumnWidth = 75.29666666666667
$ws.Columns("J").ColumnWidth = 102.66666666666667
$ws.Columns("K").ColumnWidth = 54.79666666666667
$ws.Columns("L").ColumnWidth = 94.54666666666667
$ws.Columns("M").ColumnWidth = 31.416666666666668
$ws.Columns("N").ColumnWidth = 62.666666666666664
$ws.Columns("O").ColumnWidth = 27.416666666666668
$ws.Columns("P").ColumnWidth = 120.41666666666667
$ws.Columns("Q").ColumnWidth = 110.54666666666667
$ws.Columns("R").ColumnWidth = 176.16666666666666
$ws.Columns("S").ColumnWidth = 114.04666666666667
$ws.Columns("T").ColumnWidth = 103.91666666666667
$ws.Columns("U").ColumnWidth = 52.416666666666664
$ws.Columns("V").ColumnWidth = 98.04666666666667
$ws.Columns("W").ColumnWidth = 50.29666666666667
$ws.Columns("X").ColumnWidth = 167.41666666666666
$ws.Columns("Y").ColumnWidth = 53.166666666666664
$ws.Columns("Z").ColumnWidth = 33.166666666666664
$ws.Columns("AA").ColumnWidth = 56.416666666666664
